# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
#
# For a handful of row-groups in the "Sweden 1div Norra" sheet, the match
# records (columns B..AB) were cyclically re-ordered by one position
# within the group (column A, the sequential row id, stays put). This
# reproduces that re-ordering using the Excel object model: for each
# group we snapshot the B:AB values of every row first (so later writes
# never clobber data we still need to read), then write each row the
# values that previously belonged to the row "above" it in the group,
# wrapping the first row around to take the last row's original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row groups, in on-sheet order, whose B:AB payloads get rotated by one.
$groups = @(
    @(24, 25),
    @(45, 46),
    @(74, 75, 76, 77),
    @(133, 134),
    @(146, 147),
    @(211, 212),
    @(222, 223)
)

foreach ($group in $groups) {
    # Snapshot the current B:AB contents of every row in this group first.
    $snapshots = @{}
    foreach ($row in $group) {
        $snapshots[$row] = $ws.Range("B$row`:AB$row").Value2
    }

    # Write back: row[i] receives the snapshot that belonged to row[i-1],
    # with the very first row in the group wrapping around to the last.
    $n = $group.Count
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcRow = $group[($i - 1 + $n) % $n]
        $ws.Range("B$destRow`:AB$destRow").Value2 = $snapshots[$srcRow]
    }
}
